$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style from existing header cell (H1) to new header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data rows for columns I (I0) and J (IF)
$data = @(
    @(6, 7),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(6, 7),
    @(6, 7),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(8, 8),
    @(6, 6),
    @(3, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
